# Fixed link in quick start guide to tutorial.
# Merge previously-split runs back into single runs (no textual change,
# just collapsing multiple <a:r> elements that shared identical formatting).
#
# Trick: PowerPoint's TextRange.Text setter merges all runs covered by the
# range into a single run (using the formatting of the first run) *only*
# when the assigned string differs from the current text. Since the target
# text is identical to the concatenation of the existing runs, we first set
# a unique placeholder value (forcing the merge/rewrite) and then set the
# final, real text back over the same range.

$p = $ppt.ActivePresentation

function Merge-Text($range, [string]$finalText) {
    $placeholder = "__merge_marker__" + [guid]::NewGuid().ToString("N")
    $range.Text = $placeholder
    $range2 = $range.Characters(1, $placeholder.Length)
    $range2.Text = $finalText
}

# --- Slide 3: "Text Placeholder 2" ---
# "What is this " + "quick " + "start guide about?" -> "What is this quick start guide about?"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
Merge-Text $tr3 "What is this quick start guide about?"

# --- Slide 26: "Title 1" ---
# "The .in " + "Suffix" -> "The .in Suffix"
$s26 = $p.Slides.Item(26)
$tr26 = $s26.Shapes.Item(1).TextFrame.TextRange
Merge-Text $tr26 "The .in Suffix"

# --- Slide 27: "Content Placeholder 2" ---
# Paragraph 1: "Now build the libraries" + ":" -> "Now build the libraries:"
# Paragraph 6: "And install " + "them:" -> "And install them:"
$s27 = $p.Slides.Item(27)
$tr27 = $s27.Shapes.Item(5).TextFrame.TextRange
Merge-Text $tr27.Paragraphs(1) "Now build the libraries:"
Merge-Text $tr27.Paragraphs(6) "And install them:"

# --- Slide 28: "TextBox 6" ---
# "If this was not clear enough " + "or " + "you would like to know more, have a look at the corresponding "
# -> "If this was not clear enough or you would like to know more, have a look at the corresponding "
# (leave the following hyperlink run "BASIS Tutorial" and trailing run untouched)
$s28 = $p.Slides.Item(28)
$tr28 = $s28.Shapes.Item(6).TextFrame.TextRange
$target28 = "If this was not clear enough or you would like to know more, have a look at the corresponding "
$sub28 = $tr28.Characters(1, $target28.Length)
Merge-Text $sub28 $target28
